$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 535
$ws.Range("C2").Value = 5504

$ws.Range("B3").Value = 0.34762832999999999
$ws.Range("C3").Value = 3.5763482780000002

$ws.Range("B3").Select()
